$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row B11: 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row B12: 63 -> 105
$ws.Range("B12").Value = 105

# Update Corr/total marks text E12: "62/84" -> "105/140"
$ws.Range("E12").Value = "105/140"
